$p = $ppt.ActivePresentation
$s = $p.Slides.Item(17)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# Model A - Mikro: 65% -> 72%
$para = $tr.Paragraphs(2, 1)
$found = $para.Find("65%", 0, $false, $false)
$found.Text = "72%"

# Model A - Makro: 48% -> 49%
$para = $tr.Paragraphs(3, 1)
$found = $para.Find("48%", 0, $false, $false)
$found.Text = "49%"

# Model B - Mikro: 72% -> 65%
$para = $tr.Paragraphs(5, 1)
$found = $para.Find("72%", 0, $false, $false)
$found.Text = "65%"

# Model B - Makro: 49% -> 48%
$para = $tr.Paragraphs(6, 1)
$found = $para.Find("49%", 0, $false, $false)
$found.Text = "48%"
